$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.944.38"
$ws.Range("E2").Value = "  +1.72%  "

$ws.Range("D3").Value = "1.778.09"
$ws.Range("E3").Value = "  +1.70%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.88%  "

$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4539"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.92%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3592"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.59%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07490"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.18%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.01"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.03%  "

$ws.Range("E11").Value = "  +1.59%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.02%  "

$ws.Range("E13").Value = "  +1.12%  "

$ws.Range("E14").Value = "  +0.88%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.212"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.77%  "

$ws.Range("D16").Value = "1.780.53"
$ws.Range("E16").Value = "  +1.73%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.77"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.10%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001063"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.48%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06431"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.38%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.04%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.60%  "

$ws.Range("E22").Value = "  +0.16%  "

$ws.Range("D23").Value = "27.972.24"
$ws.Range("E23").Value = "  +1.60%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.93%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.089"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.60%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.02%  "

$ws.Range("E27").Value = "  -0.55%  "

$ws.Range("D28").Value = "1.983.59"
$ws.Range("E28").Value = "  +1.78%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.227"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.62%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.99"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.33%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.117"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.26%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09215"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.91%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.679"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.28%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.561"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.25%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.91"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02303"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.94%  "

$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06185"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.92%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2094"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.08%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6346"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.983"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.18%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.189"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.387"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.71%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.918"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.36%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.37"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.69%  "

$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5921"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.73%  "

$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.733"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.52%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.90%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.963"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.07%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06932"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.57%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.139"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.16%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.97"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.08%  "
